$d = $word.ActiveDocument

function New-Run($range, $textB64, $italic) {
    $bytes = [Convert]::FromBase64String($textB64)
    $text = [System.Text.Encoding]::UTF8.GetString($bytes)
    $range.Collapse(0)
    $range.InsertAfter($text)
    if ($italic) {
        $range.Font.Italic = $true
    }
    $range.Collapse(0)
}

# Anchor: last paragraph in the main body ("BILAGA 1 - Fridlysta arter")
$paras = $d.Paragraphs
$anchor = $paras.Item($paras.Count)
$insertAfterRange = $anchor.Range

# --- New paragraph 0 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Heading 1"
$pr = $newp.Range
New-Run $pr "S27DpHJvdCDigJMgZWtvbG9naSBzYW10IGtyYXYgcMOlIGxpdnNtaWxqw7Zu" $false
$insertAfterRange = $newp.Range

# --- New paragraph 1 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Normal"
$pr = $newp.Range
New-Run $pr "S27DpHJvdCDDpHIgZnJpZGx5c3QgZW5saWd0IDggb2NoIDE1IMKnwqcgYXJ0c2t5ZGRzZsO2cm9yZG5pbmdlbiBvY2gga2xhc3NhZCBzb20gc8OlcmJhciAoVlUpIGVubGlndCByw7ZkbGlzdGFuIDIwMjAuIEtuw6Ryb3Qgw6RyIGJlcm9lbmRlIGF2IGjDtmcgb2NoIGrDpG1uIGx1ZnRmdWt0aWdoZXQgaSBnYW1sYSwgb3N0w7ZyZGEgc2tvZ3NtaWxqw7ZlciBvY2ggw6RyIGvDpG5zbGlnIGbDtnIgc25hYmJhIGbDtnLDpG5kcmluZ2FyIGF2IGxqdXMtL3ZpbmRmw7ZyaMOlbGxhbmRlbiBlbGxlciB1dHRvcmtuaW5nLiBQw6UgZ3J1bmQgYXYgZXR0IGFsbHRmw7ZyIGludGVuc2l2dCBza29nc2JydWsgaGFyIGRlbiBtaW5za2F0IG1lZCA0MCAoMjUtNTApICUgdW5kZXIgZGUgc2VuYXN0ZSA2MCDDpXJlbiBvY2ggaSBmcmFtdGlkZW4gYmVkw7ZtcyBtaW5za25pbmdzdGFrdGVuIHVwcGfDpSB0aWxsIDMwICgyMC00MCkgJS4gVGlsbCBmw7ZsamQgYXYgYXR0IGFydGVuIGhhciBlbiBkb2t1bWVudGVyYXQgaMO2Z3JlIG1pbnNrbmluZ3N0YWt0IGlmw7ZyaMOlbGxhbmRlIHRpbGwgc2luIGdlbmVyYXRpb25zdGlkIMOkbiB2YWQgc29tIHRpZGlnYXJlIHZhcml0IGvDpG50IChkYXRhIGZyw6VuIFJpa3Nza29nc3RheGVyaW5nZW4pIGjDtmpkZXMgZGVuIHRpbGwgaG90a2F0ZWdvcmkgc8OlcmJhciAoVlUpIGkgcsO2ZGxpc3RhbiAyMDIwIChBcnRkYXRhYmFua2VuLCAyMDIxKS4=" $false
$insertAfterRange = $newp.Range

# --- New paragraph 2 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Normal"
$pr = $newp.Range
New-Run $pr "U2FtdWVsIEpvaG5zb25zIGRva3RvcnNhdmhhbmRsaW5nIA==" $false
New-Run $pr "4oCcUmV0ZW50aW9uIEZvcmVzdHJ5IGFzIGEgQ29uc2VydmF0aW9uIE1lYXN1cmUgZm9yIEJvcmVhbCBGb3Jlc3QgR3JvdW5kIFZlZ2V0YXRpb27igJw=" $true
New-Run $pr "IChTTFUsIFVwcHNhbGEgMjAxNCkgdmlzYXIgYXR0IGRldCBrcsOkdnMgdsOkbCB0aWxsdGFnbmEgc2t5ZGRzem9uZXIgZsO2ciBhdHQga27DpHJvdGVucyB2w6R4dHBsYXRzZXIgaW50ZSBza2EgdGEgc2thZGEgYXYgc2tvZ3NicnVrc8OldGfDpHJkZXIgaSBpbnRpbGxpZ2dhbmRlIG9tcsOlZGVuOiA=" $false
New-Run $pr "4oCcU3R1ZHkgSUlJIHNob3dzIHRoYXQgcmV0ZW50aW9uIHBhdGNoZXMgc21hbGxlciB0aGFuIDAuNSBoYSBkbyBub3QgbGlmZWJvYXQgdGhlIHNlbnNpdGl2ZSBmb3Jlc3QgaGVyYiBHLiByZXBlbnMsIGEgc3BlY2llcyB0aGF0IGRlcGVuZCBvbiBzdGFibGUgbWljcm9jbGltYXRpYyBjb25kaXRpb25zIHR5cGljYWwgZm9yIGludGFjdCBmb3Jlc3Qgc3RhbmRzLuKAnSA=" $true
New-Run $pr "VmlkYXJlIA==" $false
New-Run $pr "4oCcTW9yZSBzZW5zaXRpdmUgZm9yZXN0IHNwZWNpZXMgYXJlIG5vdCBsaWZlYm9hdGVkIGluIHJldGVudGlvbiBwYXRjaGVzIHJhbmdpbmcgZnJvbSAwLjA1IHRvIDAuNSBoYSAoUGFwZXJzIElJICYgSUlJKS7igJ0=" $true
$insertAfterRange = $newp.Range

# --- New paragraph 3 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Normal"
$pr = $newp.Range
New-Run $pr "Sm9obnNvbnMgKDIwMTQpIHJla29tbWVuZGF0aW9uIHDDpSBtaW5zdCA1MCBtZXRlcnMgYnJlZGEgc2t5ZGRzem9uZXIgcnVudCBrbsOkcm90ZW5zIHbDpHh0cGxhdHNlciBtb3RzdmFyYXIgZW4gYXJlYWwgcMOlIDAsNzggaGVrdGFyLCB2aWxrZXQgbGlnZ2VyIGkgbGluamUgbWVkIGFuZHJhIHN0dWRpZXIgc29tIGdqb3J0cyBww6Uga8OkbnNsaWdhIHNrb2dzYXJ0ZXI6IA==" $false
New-Run $pr "4oCcSW4gc3R1ZHkgSUlJIEkgYWxzbyBzaG93IHRoYXQgdHJhbnNsb2NhdGVkIHNwZWNpbWVucyBvZiBHLiByZXBlbnMgc3Vydml2ZXMgd2VsbCBpbiBtYXR1cmUgZm9yZXN0cyBhdCBsZWFzdCA1MCBtIGZyb20gdGhlIG5lYXJlc3QgZWRnZSB0byBhbiBvcGVuIGFyZWEuIE1vcmVvdmVyLCBtZWFzdXJlcyBvZiB0ZW1wZXJhdHVyZSBhbmQgaHVtaWRpdHkgc2hvdyB0aGF0IHN1Y2ggZGlzdGFuY2VzIGZyb20gYW4gb3BlbiBhcmVhIGlzIGZhciBlbm91Z2ggdG8gb2ZmZXIgYSBtaWNyb2NsaW1hdGUgdGhhdCBpcyBtb3JlIHN0YWJsZSBjb21wYXJlZCB0byB3aGF0IHByZXNlbnQgaW4gcmV0ZW50aW9uIHBhdGNoZXMgb2YgYXJvdW5kIDAuMSBoYS4gVGhpcyBtZWFucyB0aGF0IHRoZSB2ZXJ5IGNlbnRyZSBvZiBhIGNpcmN1bGFyIHBhdGNoIHdpdGggcmFkaXVzIDUwIG0gKGVxdWFscyBhIHNpemUgb2YgMC43OCBoYSkgc2hvdWxkIG9mZmVyIGNvbmRpdGlvbnMgc2ltaWxhciB0byBpbnRlcmlvciBmb3Jlc3QgYW5kIHdvdWxkIHBlcmhhcHMgYmUgYSBzdWl0YWJsZSBoYWJpdGF0IGZvciBHLiByZXBlbnMgYW5kIHNpbWlsYXIgc3BlY2llcy4gUHJldmlvdXMgc3R1ZGllcyBmcm9tIGJvdGggTm9ydGggQW1lcmljYSBhbmQgU3dlZGVuIGhhdmUgYWxzbyBjb25jbHVkZWQgdGhhdCBwYXRjaGVzIGJldHdlZW4gMC41IGFuZCBvbmUgaGEgYXJlIHN1ZmZpY2llbnQgZm9yIHByZXNlcnZpbmcgaW50ZXJpb3IgZm9yZXN0IHZlZ2V0YXRpb24gYXMgd2VsbCBhcyBzZW5zaXRpdmUgbGljaGVucyBhbmQgYnJ5b3BoeXRlcyAoZGUgR3JhYWYgJiBSb2JlcnRzIDIwMDk7IEhhbHBlcm4gZXQgYWwuIDIwMTI7IFJ1ZG9scGhpIGV0IGFsLiAyMDE0KS7igJ0=" $true
$insertAfterRange = $newp.Range

# --- New paragraph 4 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Normal"
$pr = $newp.Range
New-Run $pr "RW4gbnlsaWdlbiBwdWJsaWNlcmFkIHZldGVuc2thcGxpZyB1cHBzYXRzIGF2IEtvZWxtZWlqZXIgbS5mbC4gKDIwMjIpIGlua2x1ZGVyYXIgb3JraWTDqW4ga27DpHJvdHMgc2t5ZGRzYmVob3YuIEkgdXBwc2F0c2VuIGJlcsO2cnMgcHJvYmxlbWV0IG1lZCB1dHRvcmtuaW5nIGbDtnIgdsOkeHRlciwgYmwuYS4gZsO2ciBrbsOkcm90LCBldHQgcHJvYmxlbSBzb20gYmxpdml0IGFjY2VudHVlcmF0IHDDpSBncnVuZCBhdiBkZW4gcMOlZ8OlZW5kZSBrbGltYXRmw7Zyw6RuZHJpbmdlbiBvY2ggdG9ycmEgc29tcmFyLCB0LmV4LiBkZW4gZXhjZXB0aW9uZWxsdCB0b3JyYSBzb21tYXJlbiAyMDE4LiBJIHVwcHNhdHNlbiB1bmRlcnPDtmtzIG9tcsOlZGVuIG1lZCB0cmUgb2xpa2EgYXZzdMOlbmQgZnLDpW4ga2FsaHlnZ2Vza2FudCBtZWQgYXZzZWVuZGUgcMOlIHNreWRkIGJsLmEuIGbDtnIga27DpHJvdC4gRGV0IGbDtnJzdGEgb21yw6VkZXQgaGFyIGF2c3TDpW5kIHVwcCB0aWxsIDIwIG0gZnLDpW4gaHlnZ2Vza2FudCAoU3Ryb25nIGVkZ2UgZWZmZWN0KSwgZGV0IGFuZHJhIDIwIOKAkyA0MCBtIGZyw6VuIGh5Z2dlc2thbnQgKFdlYWsgZWRnZSBlZmZlY3QpIG9jaCBkZXQgdHJlZGplIGF2c2VyIHN0w7ZycmUgYXZzdMOlbmQgZnLDpW4gaHlnZ2Vza2FudCwgZMOkciBrYW50ZWZmZWt0ZW4gYW5zZXMgdmFyYSBmw7Zyc3VtYmFyIChJbnRlcmlvcikuIEV0dCByZXN1bHRhdCB2YXIgYXR0IG1hbiBmYW5uIHN0b3IgZWxsZXIgbXlja2V0IHN0b3IgdXR0b3JrbmluZ3NlZmZla3QgcMOlIGvDpG5zbGlnYSBvY2ggcsO2ZGxpc3RhZGUgc2tvZ3NhcnRlciB2aWQgZGUga29ydGFyZSBhdnN0w6VuZGVuIHRpbGwgaHlnZ2Vza2FudCwgbWVkYW4gZWZmZWt0IGF2IHV0dG9ya25pbmcgaW50ZSBrb25zdGF0ZXJhZGVzIHDDpSBzdMO2cnJlIGF2c3TDpW5kIChJbnRlcmlvcikuIEbDtnIgb3JraWTDqW4ga27DpHJvdCBmYW5uIG1hbiBlbiByaWsgZsO2cmVrb21zdCAodXBwIHRpbGwgMCwwNiBkbTIvbTIpIHDDpSBzdG9ydCBhdnN0w6VuZCBmcsOlbiBoeWdnZXNrYW50IChJbnRlcmlvciksIG1lZGFuIGbDtnJla29tc3RlbiB2YXIgbGl0ZW4gZWxsZXIgbsOkcm1hc3QgZsO2cnN1bWJhciBpIGRlIG9tcsOlZGVuIHNvbSBrbGFzc2lmaWNlcmFkZXMgc29tIFdlYWsgZWRnZSBlZmZlY3QgcmVzcGVrdGl2ZSBTdHJvbmcgZWRnZSBlZmZlY3QuIEFyYmV0ZXQgcMOlcGVrYXIgYXR0IGRlIGFsbHQgb2Z0YXJlIGbDtnJla29tbWFuZGUgdG9ycmEgc29tcmFybmEgZ2VyIHl0dGVybGlnYXJlIHNrw6RsIGF0dCB1dMO2a2Egc2t5ZGRzYXZzdMOlbmRldCBmcsOlbiBoeWdnZW4gdGlsbCBkZW4gZnVrdGtyw6R2YW5kZSBhcnRlbiBrbsOkcm90IChLb2VsbWVpamVyIG0uZmwuLCAyMDIyKS4=" $false
$insertAfterRange = $newp.Range

# --- New paragraph 5 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Normal"
$pr = $newp.Range
New-Run $pr "w4R2ZW4gU2tvZ3NzdHlyZWxzZW5zIGVnZW4gdsOkZ2xlZG5pbmcgZsO2ciBow6Ruc3luIHRpbGwga27DpHJvdCBsaWdnZXIgaSBsaW5qZSBtZWQgb3ZhbnN0w6VlbmRlIGZvcnNrbmluZ3NzdHVkaWVyLiBBdiB2w6RnbGVkbmluZ2VuIGZyYW1nw6VyIGRldCBhdHQgZsO2ciBtZWQgaMO2ZyBzYW5ub2xpa2hldCBrdW5uYSBiZXZhcmEgYmVmaW50bGlnYSBmw7ZyZWtvbXN0ZXIga3LDpHZzIHJlbGF0aXZ0IHN0b3JhIGF2c8OkdHRuaW5nYXIgYXYgdXBwdnV4ZW4gc2tvZyBtZWQgc2x1dGV0IG9jaCByZWxhdGl2dCB0w6R0dCBrcm9uc2tpa3QuIFNvbSByaWt0bGluamUga2FuIGtyw6R2YXMgZXR0IGF2c3TDpW5kIHDDpSA1MCBtZXRlciBpbiBmcsOlbiBicnluZXQgZsO2ciBhdHQgdmlkbWFrdGjDpWxsYSBldHQgZnVuZ2VyYW5kZSBtaWtyb2tsaW1hdC4gRGV0dGEgaW5uZWLDpHIgYXR0IGZyaXN0w6VlbmRlIGjDpG5zeW5zeXRvciBmw7ZyIG3DpW5nYSBhcnRlciAoa8Okcmx2w6R4dGVyLCBsYXZhciBvY2ggbW9zc29yKSBrYW4gYmVow7Z2YSBoYSBlbiBhcmVhIMO2dmVyc3RpZ2FuZGUgMCw4IGhla3RhciAoY2lya2VseXRhIG1lZCByYWRpZW4gNTAgbWV0ZXIgPSAwLDc4IGhla3RhcikgZsO2ciBhdHQgYmliZWjDpWxsYSBsb2thbGtsaW1hdGV0LiDDhHZlbiBnYW5za2Egc23DpSBmw7Zyw6RuZHJpbmdhciBpIGZvcm0gYXYgZsO2csOkbmRyYWRlIGxqdXMtIG9jaCBmdWt0aWdoZXRzZsO2cmjDpWxsYW5kZW4sIHRpbGwgZXhlbXBlbCB0aWxsIGbDtmxqZCBhdiBnYWxscmluZywga2FuIGxlZGEgdGlsbCBhdHQgYXJ0ZW4gZsO2cnN2aW5uZXIgdGlsbCBmw7ZsamQgYXYga29ua3VycmVucyBtZWQgbWVyYSBsanVza3LDpHZhbmRlIG9jaCBzbmFiYnbDpHhhbmRlIGFydGVyIChTa29nc3N0eXJlbHNlbiwgMjAyMiku" $false
$insertAfterRange = $newp.Range

# --- New paragraph 6 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Heading 2"
$pr = $newp.Range
New-Run $pr "UmVmZXJlbnNlciAtIGtuw6Ryb3Q=" $false
$insertAfterRange = $newp.Range

# --- New paragraph 7 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Normal"
$pr = $newp.Range
New-Run $pr "ZGUgR3JhYWYgTSAmIFJvYmVydHMgTS5SLiwgMjAwOS4g" $false
New-Run $pr "U2hvcnQtdGVybSByZXNwb25zZSBvZiB0aGUgaGVyYmFjZW91cyBsYXllciB3aXRoaW4gbGVhdmUgcGF0Y2hlcyBhZnRlciBoYXJ2ZXN0LiA=" $true
New-Run $pr "Rm9yZXN0IEVjb2xvZ3kgYW5kIE1hbmFnZW1lbnQgMjU3LCAxMDE0LTEwMjU=" $false
$insertAfterRange = $newp.Range

# --- New paragraph 8 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Normal"
$pr = $newp.Range
New-Run $pr "SGFscGVybiwgQy4gQi4sIEhhbGFqLCBKLiwgRXZhbnMsIFMuIEEuLCAmIERvdmNpYWssIE0uLCAyMDEyLiA=" $false
New-Run $pr "TGV2ZWwgYW5kIHBhdHRlcm4gb2Ygb3ZlcnN0b3J5IHJldGVudGlvbiBpbnRlcmFjdCB0byBzaGFwZSBsb25nLXRlcm0gcmVzcG9uc2VzIG9mIHVuZGVyc3RvcmllcyB0byB0aW1iZXIgaGFydmVzdC4g" $true
New-Run $pr "RWNvbG9naWNhbCBBcHBsaWNhdGlvbnMsIDIyLCAyMDQ5LTIwNjQg" $false
$insertAfterRange = $newp.Range

# --- New paragraph 9 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Normal"
$pr = $newp.Range
New-Run $pr "S29lbG1laWplciwgSS4gQS4sIEVocmzDqW4sIEouLCBKw7Zuc3NvbiwgTS4sIERlIEZyZW5uZSwgUC4sIEJlcmcsIFAuLCBBbmRlcnNzb24sIEouLCBXZWlidWxsLCBILiAmIEh5bGFuZGVyLCBOLiAyMDIyLiA=" $false
New-Run $pr "SW50ZXJhY3RpdmUgZWZmZWN0cyBvZiBkcm91Z2h0IGFuZCBlZGdlIGV4cG9zdXJlIG9uIG9sZC1ncm93dGggZm9yZXN0IHVuZGVyc3Rvcnkgc3BlY2llcy4g" $true
New-Run $pr "TGFuZHNjYXBlIEVjb2xvZ3ksIDM3LCBzaWQgMTgzOS0xODUz" $false
$insertAfterRange = $newp.Range

# --- New paragraph 10 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Normal"
$pr = $newp.Range
New-Run $pr "UnVkb2xwaGksIEouLCBKw7Zuc3NvbiwgTS4gVC4sICYgR3VzdGFmc3NvbiwgTC4sIDIwMTQuIA==" $false
New-Run $pr "QmlvbG9naWNhbCBsZWdhY2llcyBidWZmZXIgbG9jYWwgc3BlY2llcyBleHRpbmN0aW9uIGFmdGVyIGxvZ2dpbmcuIA==" $true
New-Run $pr "Sm91cm5hbCBvZiBBcHBsaWVkIEVjb2xvZ3kuIDUxLCA1My02Mi4=" $false
$insertAfterRange = $newp.Range

# --- New paragraph 11 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Normal"
$pr = $newp.Range
New-Run $pr "U2tvZ3NzdHlyZWxzZW4sIDIwMjIuIA==" $false
New-Run $pr "VsOkZ2xlZG5pbmcgZsO2ciBow6Ruc3luIHRpbGwga27DpHJvdC4g" $true
New-Run $pr "aHR0cHM6Ly93d3cuc2tvZ3NzdHlyZWxzZW4uc2UvbGFnLW9jaC10aWxsc3luL2FydHNreWRkL3ZhZ2xlZG5pbmdhci1vY2gta3Vuc2thcHNzdG9kLWFydHNreWRkL3ZhZ2xlZG5pbmctZm9yLWhhbnN5bi10aWxsLWtuYXJvdC8=" $false
$insertAfterRange = $newp.Range

# --- New paragraph 12 ---
$newRange = $insertAfterRange.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($n)
$newp.Style = "Normal"
$pr = $newp.Range
New-Run $pr "U0xVIEFydGRhdGFiYW5rZW4sIDIwMjEuIA==" $false
New-Run $pr "QXJ0ZmFrdGFibGFkLiBOYXR1cnbDpXJkIOKAkyBhcnRmYWt0YS4g" $true
New-Run $pr "U0xVIEFydGRhdGFiYW5rZW4sIFVwcHNhbGEg" $false
$insertAfterRange = $newp.Range

# --- Update date in the first-page header ---
$sec = $d.Sections.Item(1)
$headers = $sec.Headers
$firstPageHeader = $headers.Item(2)
$firstPageHeader.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null

Write-Output "Done. Paragraph count now: $($d.Paragraphs.Count)"
